$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.545.75"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "1.844.12"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3858"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07891"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9979"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.966"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "1.842.04"
$ws.Range("E14").Value = "  -3.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.127"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06670"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001034"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.006"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").Value = "27.546.95"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.315"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").Value = "2.064.20"
$ws.Range("E27").Value = "  -3.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.118"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.408"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9759"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09409"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.601"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.305"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.329"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06032"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02227"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.302"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.180"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5903"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1863"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.240"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5570"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.910"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06697"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.81%  "
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.007"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.14%  "
